# 2019 TN crash data: a new "Unknown" county row was extracted/merged into
# the cleaned dataset and inserted in its correct alphabetically-sorted
# position between "Union" (row 88) and "VanBuren" (former row 89, now 90).
# All the rows below shift down by one; no other existing data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push "VanBuren" ... "Wilson" (old rows 89-96) down to make room.
$ws.Rows.Item(89).Insert()

# New row: Unknown county — 0 fatalities, 0 injuries, 1 total crash, and no
# CBSA / metro-micro info (those three columns stay blank, matching the
# other sparsely-populated counties already in the sheet).
$ws.Range("A89").Value = "Unknown"

# fatalities/injuries/total_crashes are stored as text in this sheet (see
# column B/C/D elsewhere, e.g. "4", "1,911"), so force text formatting
# before writing the numeric-looking strings.
$numCells = $ws.Range("B89:D89")
$numCells.NumberFormat = "@"
$ws.Range("B89").Value = "0"
$ws.Range("C89").Value = "0"
$ws.Range("D89").Value = "1"
